$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A from slash format (DD/MM/YYYY) to dash
# format (DD-MM-YYYY) for rows 3-21. A leading apostrophe is used so Excel
# stores the values as plain text instead of re-interpreting the
# dash-separated strings as date serials.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $ws.Cells.Item($r, 1).Formula = "'" + $dates[$r]
}

# Update the attendance count values that changed
$ws.Cells.Item(3, 4).Value = 1   # D3: Total Attendance Count
$ws.Cells.Item(3, 7).Value = 1   # G3: Invalid

$ws.Cells.Item(4, 4).Value = 1   # D4: Total Attendance Count
$ws.Cells.Item(4, 5).Value = 1   # E4: Real
$ws.Cells.Item(4, 8).Value = 0   # H4: Absent
